# pedidos.xlsx edit script
# - removes the "Cantidad" column (old column B)
# - inserts a new "Precio total" column before "Precio cancelado"
# - merges product + quantity info into the "Tipo de producto" column (col A)
# - fixes/normalizes the exported date values (kept as raw numeric serials)
# - fixes hyperlinks so they point at the "Correo electrónico" column after the
#   column shift

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Remove existing hyperlinks (they will be re-created after the columns
#    are rearranged, since their target cell moves from F to E).
# ---------------------------------------------------------------------------
$ws.Range("F2").Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 2) Remove the "Cantidad" column (old column B) and insert a new blank
#    column for "Precio total" right before "Precio cancelado" (old column I,
#    which after the deletion above is column H).
# ---------------------------------------------------------------------------
$ws.Columns("B:B").Delete()
$ws.Columns("H:H").Insert()

# The new column inherits the text ("@") format of the column to its left
# (old "Estado" column), which would turn the new numeric prices into text.
# Reset that formatting so the "Precio total" values are stored as real
# numbers.
$ws.Range("H1:H6").ClearFormats()

# ---------------------------------------------------------------------------
# 3) Header row
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Tipo de producto"
$ws.Range("B1").Value = "Fecha de solicitud"
$ws.Range("C1").Value = "Fecha de retiro estimado"
$ws.Range("D1").Value = "Nombre"
$ws.Range("E1").Value = "Correo electrónico"
$ws.Range("F1").Value = "Número telefónico"
$ws.Range("G1").Value = "Estado"
$ws.Range("H1").NumberFormat = "0"
$ws.Range("H1").Value = "Precio total"
$ws.Range("I1").Value = "Precio cancelado"
$ws.Range("J1").Value = "Porcentaje descuento"

# ---------------------------------------------------------------------------
# 4) Data rows: column A now stores a combined "qty-product" description
#    instead of splitting quantity into its own column. Dates (B,C) keep
#    their original numeric date serials. A new "Precio total" value is
#    added in column H.
# ---------------------------------------------------------------------------

# Row 2
$ws.Range("A2").Value = "1-torta,2-queques"
$ws.Range("B2").Value = 43743
$ws.Range("C2").Value = 43744
$ws.Range("D2").Value = "alberto hurtado"
$ws.Range("E2").Value = "albhurt@gmail.com"
$ws.Range("F2").Value = "9-48485930"
$ws.Range("G2").Value = "pendiente"
$ws.Range("H2").Value = 1000

# Row 3
$ws.Range("A3").Value = "1-kuchen manzana"
$ws.Range("B3").Value = 43744
$ws.Range("C3").Value = 43745
$ws.Range("D3").Value = "carmen ramirez"
$ws.Range("E3").Value = "cramirez@gmail.com"
$ws.Range("F3").Value = "9-47382938"
$ws.Range("G3").Value = "pendiente"
$ws.Range("H3").Value = 2000

# Row 4
$ws.Range("A4").Value = "1-trozo queque"
$ws.Range("B4").Value = 43745
$ws.Range("C4").Value = 43745
$ws.Range("D4").Value = "marcelo muñoz"
$ws.Range("E4").Value = "mmuñoz@gmail.com"
$ws.Range("F4").Value = "9-62844758"
$ws.Range("G4").Value = "pendiente"
$ws.Range("H4").Value = 3000

# Row 5
$ws.Range("A5").Value = "1-pie de limon"
$ws.Range("B5").Value = 43768
$ws.Range("C5").Value = 43769
$ws.Range("D5").Value = "pedro casillas"
$ws.Range("E5").Value = "pedrosillas@gmai.com"
$ws.Range("F5").Value = "9-37284547"
$ws.Range("G5").Value = "pendiente"
$ws.Range("H5").Value = 4900

# Row 6
$ws.Range("A6").Value = "1-tartaleta durazno"
$ws.Range("B6").Value = 43767
$ws.Range("C6").Value = 43799
$ws.Range("D6").Value = "lionel mauro"
$ws.Range("E6").Value = "limau@gmail.com"
$ws.Range("F6").Value = "9-38274654"
$ws.Range("G6").Value = "pendiente"
$ws.Range("H6").Value = 5000

# ---------------------------------------------------------------------------
# 5) Re-create the mailto hyperlinks on column E (Correo electrónico).
#    Adding a hyperlink forces Excel's default "Hyperlink" look (blue,
#    underlined); restore the plain look the sheet originally used for
#    these mailto links (normal color, no underline).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:albhurt@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:cramirez@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:mmuñoz@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E5"), "mailto:pedrosillas@gmai.com")
$ws.Hyperlinks.Add($ws.Range("E6"), "mailto:limau@gmail.com")

$ws.Range("E2:E6").Font.Underline = -4142
$ws.Range("E2:E6").Font.ThemeColor = 1

# ---------------------------------------------------------------------------
# 6) Selection / active cell cosmetics
# ---------------------------------------------------------------------------
$ws.Range("H6").Select()
